$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit performs a rotation of the species-observation data (columns
# A, B, E, F, G, H, Q, R plus the special annotation cells M/AC for rows
# 10/11) across rows 6, 8, 9, 10, 11 and 13, while all other
# location/date/observer columns remain fixed per row.
#
# Cycle: 6 -> 10 -> 11 -> 8 -> 6   (row N receives what row X used to hold)
# Swap:  9 <-> 13

# Row 6 receives what row 10 used to hold
$ws.Range("A6").Value = 112038121
$ws.Range("B6").Value = 89423
$ws.Range("E6").Value = 5432
$ws.Range("F6").Value = "Granticka"
$ws.Range("G6").Value = "Porodaedalea chrysoloma"
$ws.Range("H6").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q6").Value = 562313.4042944111
$ws.Range("R6").Value = 7307732.191837165

# Row 8 receives what row 6 used to hold
$ws.Range("A8").Value = 112038104
$ws.Range("B8").Value = 89401
$ws.Range("E8").Value = 1108
$ws.Range("F8").Value = "Harticka"
$ws.Range("G8").Value = "Pelloporus leporinus"
$ws.Range("H8").Value = "(Fr.) Krieglst."
$ws.Range("Q8").Value = 562291.6665034146
$ws.Range("R8").Value = 7307714.531584828

# Row 9 receives what row 13 used to hold
$ws.Range("A9").Value = 112038198
$ws.Range("Q9").Value = 562332.0260024283
$ws.Range("R9").Value = 7307761.242099251

# Row 10 receives what row 11 used to hold (including the special
# annotation cells: M10 "äldre spår" and AC10 "Skalad gran")
$ws.Range("A10").Value = 112038209
$ws.Range("B10").Value = 56398
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = "Tretåig hackspett"
$ws.Range("G10").Value = "Picoides tridactylus"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("J10").ClearContents()
$ws.Range("L10").Value = ""
$ws.Range("M10").Value = "äldre spår"
$ws.Range("Q10").Value = 562326.1318858962
$ws.Range("R10").Value = 7307824.949321065
$ws.Range("AC10").Value = "Skalad gran"
$ws.Range("AF10").ClearContents()

# Row 11 receives what row 8 used to hold (standard empty annotation
# cells, losing the special M/AC cells)
$ws.Range("A11").Value = 112038151
$ws.Range("B11").Value = 89405
$ws.Range("E11").Value = 1202
$ws.Range("F11").Value = "Ullticka"
$ws.Range("G11").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H11").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("J11").Value = ""
$ws.Range("L11").ClearContents()
$ws.Range("M11").ClearContents()
$ws.Range("Q11").Value = 562332.0260024283
$ws.Range("R11").Value = 7307761.242099251
$ws.Range("AC11").ClearContents()
$ws.Range("AF11").Value = ""

# Row 13 receives what row 9 used to hold
$ws.Range("A13").Value = 112038107
$ws.Range("Q13").Value = 562291.6665034146
$ws.Range("R13").Value = 7307714.531584828
